$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 96

# Copy the date/time style (used by the rest of column A) down onto the new row
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats

# A: date (serial number)
$ws.Cells.Item($row, 1).Value = 45461.2916666667

# B: volume
$ws.Cells.Item($row, 2).Value = 5000

# C: high
$ws.Cells.Item($row, 3).Value = 0.704999983310699

# D: low
$ws.Cells.Item($row, 4).Value = 0.699999988079071

# E: open
$ws.Cells.Item($row, 5).Value = 0.704999983310699

# F: close
$ws.Cells.Item($row, 6).Value = 0.699999988079071

# G: adj_close - stored as text (matches the rest of the column) even though it
# looks like a number, so force text entry with a leading apostrophe, then
# clear the resulting quote-prefix cell style back to Normal.
$ws.Cells.Item($row, 7).Value = "'0.699999988079071"
$ws.Cells.Item($row, 7).Style = "Normal"

# H: ticker
$ws.Cells.Item($row, 8).Value = "BWZ.MI"
